# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# listing with the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @(newPriceTextOrNull, newVolumeTextOrNull)
$updates = @{
    2 = @('35.264.67', '  -0.27%  ')
    3 = @('1.909.58', '  +0.10%  ')
    4 = @($null, '  +0.04%  ')
    5 = @('0.722', '  +8.24%  ')
    6 = @('255.96', '  +3.85%  ')
    7 = @($null, '  +0.07%  ')
    8 = @('40.69', '  -1.87%  ')
    9 = @('0.373', '  +6.82%  ')
    10 = @('52.89', '  +0.04%  ')
    11 = @('0.0760', '  +5.31%  ')
    12 = @($null, '  -0.64%  ')
    13 = @('2.186.86', '  +0.31%  ')
    14 = @('12.87', '  +6.37%  ')
    15 = @('0.729', '  +4.18%  ')
    16 = @('4.96', '  +1.82%  ')
    17 = @('1.866.62', '  -2.06%  ')
    18 = @('35.275.39', '  -0.21%  ')
    19 = @('74.65', '  +3.24%  ')
    20 = @($null, '  +2.52%  ')
    21 = @('243.83', '  +1.57%  ')
    22 = @($null, '  +4.18%  ')
    23 = @($null, '  +5.64%  ')
    24 = @($null, '  +0.07%  ')
    25 = @('2.43', '  +5.80%  ')
    26 = @('2.43', '  +4.49%  ')
    27 = @('166.31', '  -2.52%  ')
    28 = @('8.69', '  +2.95%  ')
    29 = @($null, '  +1.60%  ')
    30 = @($null, '  +3.72%  ')
    31 = @('4.127.00', '  +19.41%  ')
    32 = @('4.38', '  +5.66%  ')
    33 = @('2.00', '  +14.42%  ')
    34 = @('1.64', '  +22.09%  ')
    35 = @('0.0587', '  +3.66%  ')
    36 = @($null, '  +2.99%  ')
    37 = @($null, '  +0.02%  ')
    38 = @('0.913', '  -2.01%  ')
    39 = @($null, '  -0.39%  ')
    40 = @('17.28', '  +5.78%  ')
    41 = @('0.0219', '  +4.63%  ')
    42 = @('96.56', '  +7.11%  ')
    43 = @($null, '  +1.06%  ')
    44 = @('0.0653', '  +3.14%  ')
    45 = @('1.337.09', '  -0.33%  ')
    46 = @($null, '  +1.28%  ')
    47 = @($null, '  +1.00%  ')
    48 = @($null, '  +2.68%  ')
    49 = @($null, '  -0.73%  ')
    50 = @('44.81', '  -6.14%  ')
    51 = @('0.0751', '  +6.20%  ')
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $priceText = $pair[0]
    $volumeText = $pair[1]

    if ($priceText -ne $null) {
        $priceCell = $ws.Range("D$row")
        # The price column stores plain display text (e.g. "35.264.67"
        # or "0.722"). Values that parse as a plain number would
        # otherwise be auto-converted to a numeric cell by Excel, so
        # force text entry with a leading apostrophe and then drop the
        # resulting quote-prefix formatting to restore the original
        # (unstyled) cell style.
        if ($priceText -match '^[+-]?[0-9]*\.?[0-9]+$') {
            $priceCell.Value = "'" + $priceText
            $priceCell.Style = 'Normal'
        } else {
            $priceCell.Value = $priceText
        }
    }

    if ($volumeText -ne $null) {
        $ws.Range("E$row").Value = $volumeText
    }
}
